$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the "完成情况" (C) and "备注" (D) columns for the member rows (3-8)
# with the shared "已完成" / progress-note text.
$doneText = "已完成"
$noteText = "已完成用例图设计，下一阶段开始编写用例描述"

$ws.Range("C3").Value = $doneText
$ws.Range("D3").Value = $noteText

$ws.Range("C4").Value = $doneText
$ws.Range("D4").Value = $noteText

$ws.Range("C5").Value = $doneText
$ws.Range("D5").Value = $noteText

$ws.Range("C6").Value = $doneText
$ws.Range("D6").Value = $noteText

$ws.Range("C7").Value = $doneText
$ws.Range("D7").Value = $noteText

$ws.Range("C8").Value = $doneText
$ws.Range("D8").Value = $noteText

# Insert a new blank row above the closing "总结：" block (currently row 18),
# pushing it down to rows 19:20, and give the new row the same formatting
# as the rows above it (the plain bordered row style).
$ws.Rows.Item(18).Insert() | Out-Null
$ws.Range("A17:D17").Copy()
$ws.Range("A18:D18").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Leave the selection on D3, matching the saved state of the workbook.
$ws.Range("D3").Select() | Out-Null
